$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Person sheet: drop "current_address", remove the gender dropdown
# validation, and reorder the header row to:
# id, primary_email, name, nick, position, birth_date, gender, avatar,
# has_employment_history, aliases, description, image
# ---------------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")

# Remove the gender list data validation (was anchored on column C).
$wsPerson.Range("C2:C1048576").Validation.Delete()

# Remove the "current_address" column (column D) entirely.
$wsPerson.Columns.Item(4).Delete()

# Rewrite the header row in the new column order.
$wsPerson.Range("A1").Value = "id"
$wsPerson.Range("B1").Value = "primary_email"
$wsPerson.Range("C1").Value = "name"
$wsPerson.Range("D1").Value = "nick"
$wsPerson.Range("E1").Value = "position"
$wsPerson.Range("F1").Value = "birth_date"
$wsPerson.Range("G1").Value = "gender"
$wsPerson.Range("H1").Value = "avatar"
$wsPerson.Range("I1").Value = "has_employment_history"
$wsPerson.Range("J1").Value = "aliases"
$wsPerson.Range("K1").Value = "description"
$wsPerson.Range("L1").Value = "image"

# ---------------------------------------------------------------------------
# Context sheet: rename "mission_statement" -> "purpose_statement"
# ---------------------------------------------------------------------------
$wsContext = $wb.Worksheets.Item("Context")
$wsContext.Range("B1").Value = "purpose_statement"

# ---------------------------------------------------------------------------
# Role sheet: add a new "status" column before "aliases"
# ---------------------------------------------------------------------------
$wsRole = $wb.Worksheets.Item("Role")
$wsRole.Columns.Item(6).Insert()
$wsRole.Range("F1").Value = "status"

# ---------------------------------------------------------------------------
# Membership sheet: add a new "status" column before "description"
# ---------------------------------------------------------------------------
$wsMembership = $wb.Worksheets.Item("Membership")
$wsMembership.Columns.Item(6).Insert()
$wsMembership.Range("F1").Value = "status"
